$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows right after row 420 (pushes the old rows 421..442 down
# to become 423..444, duplicating the old 441/442 "Betarraga" entries onto
# the new 443/444 in the process since the sheet only had 24 data rows for
# this block).
$ws.Rows.Item(421).Resize(2).Insert()

# --- Row 419 (Primera) gets this week's new figures -----------------------
$ws.Range("D419").Value = 45021
$ws.Range("K419").Value = 700
$ws.Range("L419").Value = 800
$ws.Range("M419").Value = 750
$ws.Range("P419").Value = 188

# --- Row 420 (Segunda) gets this week's new figures ------------------------
$ws.Range("D420").Value = 45021
$ws.Range("J420").Value = 600
$ws.Range("K420").Value = 700
$ws.Range("L420").Value = 800
$ws.Range("M420").Value = 750
$ws.Range("P420").Value = 150

# --- New row 421 (Primera) - carries what used to be row 419's data --------
$ws.Range("A421").Value = 1
$ws.Range("B421").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C421").Value = "Arica y Parinacota"
$ws.Range("D421").Value = 44414
$ws.Range("E421").Value = 15
$ws.Range("F421").Value = 100114014
$ws.Range("G421").Value = "Betarraga"
$ws.Range("H421").Value = "Sin especificar"
$ws.Range("I421").Value = "Primera"
$ws.Range("J421").Value = 700
$ws.Range("K421").Value = 500
$ws.Range("L421").Value = 550
$ws.Range("M421").Value = 525
$ws.Range("N421").Value = "`$/paquete 4 unidades"
$ws.Range("O421").Value = "Región de Arica y Parinacota"
$ws.Range("P421").Value = 131
$ws.Range("Q421").Value = 4
$ws.Range("R421").Value = "Hortaliza"

# --- New row 422 (Segunda) - carries what used to be row 420's data --------
$ws.Range("A422").Value = 1
$ws.Range("B422").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C422").Value = "Arica y Parinacota"
$ws.Range("D422").Value = 44414
$ws.Range("E422").Value = 15
$ws.Range("F422").Value = 100114014
$ws.Range("G422").Value = "Betarraga"
$ws.Range("H422").Value = "Sin especificar"
$ws.Range("I422").Value = "Segunda"
$ws.Range("J422").Value = 900
$ws.Range("K422").Value = 500
$ws.Range("L422").Value = 550
$ws.Range("M422").Value = 525
$ws.Range("N422").Value = "`$/paquete 5 unidades"
$ws.Range("O422").Value = "Región de Arica y Parinacota"
$ws.Range("P422").Value = 105
$ws.Range("Q422").Value = 5
$ws.Range("R422").Value = "Hortaliza"

# --- New row 443 (Primera) - duplicate of row 441's (pre-shift) data -------
$ws.Range("A443").Value = 1
$ws.Range("B443").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C443").Value = "Arica y Parinacota"
$ws.Range("D443").Value = 44918
$ws.Range("E443").Value = 15
$ws.Range("F443").Value = 100114014
$ws.Range("G443").Value = "Betarraga"
$ws.Range("H443").Value = "Sin especificar"
$ws.Range("I443").Value = "Primera"
$ws.Range("J443").Value = 1200
$ws.Range("K443").Value = 400
$ws.Range("L443").Value = 500
$ws.Range("M443").Value = 442
$ws.Range("N443").Value = "`$/paquete 4 unidades"
$ws.Range("O443").Value = "Región de Arica y Parinacota"
$ws.Range("P443").Value = 110
$ws.Range("Q443").Value = 4
$ws.Range("R443").Value = "Hortaliza"

# --- New row 444 (Segunda) - duplicate of row 442's (pre-shift) data -------
$ws.Range("A444").Value = 1
$ws.Range("B444").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C444").Value = "Arica y Parinacota"
$ws.Range("D444").Value = 44918
$ws.Range("E444").Value = 15
$ws.Range("F444").Value = 100114014
$ws.Range("G444").Value = "Betarraga"
$ws.Range("H444").Value = "Sin especificar"
$ws.Range("I444").Value = "Segunda"
$ws.Range("J444").Value = 1400
$ws.Range("K444").Value = 400
$ws.Range("L444").Value = 500
$ws.Range("M444").Value = 443
$ws.Range("N444").Value = "`$/paquete 5 unidades"
$ws.Range("O444").Value = "Región de Arica y Parinacota"
$ws.Range("P444").Value = 89
$ws.Range("Q444").Value = 5
$ws.Range("R444").Value = "Hortaliza"
